$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Slide 3 title: "What drives churn?" -> "What influences churn?"
# ------------------------------------------------------------------
$slide = $p.Slides.Item(3)
if ($slide.Shapes.HasTitle) {
    $slide.Shapes.Title.TextFrame.TextRange.Text = "What influences churn?"
} else {
    $slide.Shapes.Item(1).TextFrame.TextRange.Text = "What influences churn?"
}

# ------------------------------------------------------------------
# 2) Refresh the cached "datetimeFigureOut" footer date from 2/12/2025
#    to 2/13/2025 everywhere it is placed: the Slide Master and every
#    Slide Layout's Date placeholder (ppPlaceholderDate = 16).
# ------------------------------------------------------------------
$newDate = "2/13/2025"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $phType = $null
        try { $phType = $sh.PlaceholderFormat.Type } catch {}
        if ($phType -eq $ppPlaceholderDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide Master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every Slide Layout under the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
